$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain the same event data table and
# need the same updates to the "想去人数" (F) column for several rows.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 7740
    $ws.Range("F6").Value = 41
    $ws.Range("F9").Value = 5986
    $ws.Range("F17").Value = 147
}
